$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a literal (non-autoconverted) TEXT value into a cell while
# leaving its current number format / style untouched. We go through a
# temporary formula ("=<value>") and then "Paste Special -> Values" over the
# same cell; this forces Excel to store the result as a plain shared-string
# cell instead of re-parsing the text as a date/number, and it does not
# introduce any new cell style.
# ---------------------------------------------------------------------------

function Set-TextValue($rng, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)   # xlPasteValues
}

# ---------------------------------------------------------------------------
# 1) Grab the red/pink "last attendance" highlight format (fillId 2 / s=2)
#    from K6 before it is touched, so it can be stamped onto the new L-column
#    cells that need it.
# ---------------------------------------------------------------------------
$ws.Range("K6").Copy()
$highlightSource = $ws.Range("K6")

# Apply the highlighted format to the new L cells that need it (L6, L11,
# L18, L20, L25) now, while K6 still carries style s="2".
$ws.Range("L6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L11").PasteSpecial(-4122)
$ws.Range("L18").PasteSpecial(-4122)
$ws.Range("L20").PasteSpecial(-4122)
$ws.Range("L25").PasteSpecial(-4122)

# Also copy the header format (bold + border, s="1") from K1 onto L1.
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Populate the new "2025-10-21" attendance column (L).
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("L1") "2025-10-21"

Set-TextValue $ws.Range("L2") "0"
Set-TextValue $ws.Range("L3") "0"
Set-TextValue $ws.Range("L4") "0"
Set-TextValue $ws.Range("L5") "0"
Set-TextValue $ws.Range("L6") "1"
Set-TextValue $ws.Range("L7") "0"
Set-TextValue $ws.Range("L8") "0"
Set-TextValue $ws.Range("L9") "0"
Set-TextValue $ws.Range("L10") "0"
Set-TextValue $ws.Range("L11") "1"
Set-TextValue $ws.Range("L12") "0"
Set-TextValue $ws.Range("L13") "0"
Set-TextValue $ws.Range("L14") "0"
Set-TextValue $ws.Range("L15") "0"
Set-TextValue $ws.Range("L16") "0"
Set-TextValue $ws.Range("L17") "0"
Set-TextValue $ws.Range("L18") "1"
Set-TextValue $ws.Range("L19") "0"
Set-TextValue $ws.Range("L20") "1"
Set-TextValue $ws.Range("L21") "0"
Set-TextValue $ws.Range("L22") "0"
Set-TextValue $ws.Range("L23") "0"
Set-TextValue $ws.Range("L24") "0"

# L25 is stored as a genuine number (1), unlike its neighbours.
$ws.Range("L25").Value = 1

Set-TextValue $ws.Range("L26") "0"
Set-TextValue $ws.Range("L27") "0"
Set-TextValue $ws.Range("L28") "0"

# ---------------------------------------------------------------------------
# 3) The "last attendance" highlight now lives on column L, so strip the old
#    style from the K-column cells that used to carry it. The cell text
#    itself (1 / 1 / 1 / 1 / 0.5) is unchanged, only the formatting moves.
# ---------------------------------------------------------------------------
$ws.Range("K4").Style = "Normal"
$ws.Range("K6").Style = "Normal"
$ws.Range("K20").Style = "Normal"
$ws.Range("K26").Style = "Normal"
$ws.Range("K29").Style = "Normal"

# K30 additionally changes from a numeric 1 to a literal text "1".
Set-TextValue $ws.Range("K30") "1"
$ws.Range("K30").Style = "Normal"
